$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.671.87"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "1.794.33"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'308.98"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4463"
$ws.Range("E7").Value = "  +5.64%  "
$ws.Range("D8").Value = "'0.3675"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.07333"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'0.8565"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "1.815.59"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'6.609"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "'92.53"
$ws.Range("E14").Value = "  +3.44%  "
$ws.Range("D15").Value = "'0.07071"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "'5.267"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "'14.77"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").Value = "26.700.34"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").Value = "'5.138"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'10.78"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "'151.47"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'18.34"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'2.158"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("D28").Value = "'5.184"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "'117.12"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "'0.08770"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'0.7386"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "'1.149"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.439"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.903"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").Value = "'0.9998"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").Value = "'0.01952"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.05163"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'0.5308"
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("D40").Value = "'2.828"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").Value = "'6.928"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").Value = "'0.1678"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'0.5057"
$ws.Range("E43").Value = "  +6.86%  "
$ws.Range("D44").Value = "'8.388"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'10.47"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "'1.940"
$ws.Range("E46").Value = "  +3.20%  "
$ws.Range("D47").Value = "'105.01"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").Value = "'0.9996"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "'0.06285"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "'0.9140"
$ws.Range("E51").Value = "  +0.91%  "
